$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row labels (B column) for data rows 2..17, continuing the existing
# line1..line6 / extr1..extr8 naming scheme with two new "line7"/"line8"
# entries inserted before the extr* block.
$names = @("line1","line2","line3","line4","line5","line6","line7","line8","extr1","extr2","extr3","extr4","extr5","extr6","extr7","extr8")

# Target data per row: index (A), from_bus (C), to_bus (D), in_service (E)
$data = @(
    @(0, 7, 9, $true),
    @(1, 9, 8, $true),
    @(2, 8, 10, $false),
    @(3, 8, 11, $true),
    @(4, 10, 5, $true),
    @(5, 12, 8, $true),
    @(6, 14, 11, $true),
    @(7, 16, 9, $true),
    @(8, 5, 12, $true),
    @(9, 5, 9, $true),
    @(10, 10, 11, $true),
    @(11, 7, 8, $true),
    @(12, 9, 11, $true),
    @(13, 7, 11, $false),
    @(14, 5, 7, $false),
    @(15, 8, 5, $true)
)

# Before adding new rows, copy the formatting of the last existing data row
# (row 15, column A carries the bordered/bold style) down to the two new
# rows 16 and 17 so they pick up the same cell style.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
